$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '68.502.29'
$ws.Range("E2").Value = '  -1.63%  '
$ws.Range("D3").Value = '2.453.20'
$ws.Range("E3").Value = '  -1.86%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '562.18'
$ws.Range("E5").Value = '  -2.43%  '
$ws.Range("D6").Value = '163.59'
$ws.Range("E6").Value = '  -1.88%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("D9").Value = '2.451.95'
$ws.Range("E9").Value = '  -1.88%  '
$ws.Range("E10").Value = '  -5.25%  '
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("E12").Value = '  -3.96%  '
$ws.Range("E13").Value = '  -2.49%  '
$ws.Range("D14").Value = '2.905.13'
$ws.Range("E14").Value = '  -1.75%  '
$ws.Range("D15").Value = '68.408.83'
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("E16").Value = '  -3.39%  '
$ws.Range("D17").Value = '23.59'
$ws.Range("E17").Value = '  -4.70%  '
$ws.Range("D18").Value = '2.468.15'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = '10.97'
$ws.Range("E19").Value = '  -1.92%  '
$ws.Range("D20").Value = '343.26'
$ws.Range("E20").Value = '  -1.39%  '
$ws.Range("E21").Value = '  -4.36%  '
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '1.87'
$ws.Range("D25").Value = '68.09'
$ws.Range("E25").Value = '  -3.14%  '
$ws.Range("E26").Value = '  -5.40%  '
$ws.Range("D27").Value = '1.04'
$ws.Range("E27").Value = '  +3.81%  '
$ws.Range("D28").Value = '2.580.00'
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("E29").Value = '  -6.26%  '
$ws.Range("D30").Value = '0.0₃0840'
$ws.Range("E30").Value = '  -5.59%  '
$ws.Range("E31").Value = '  -6.45%  '
$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  +124.03%  '
$ws.Range("D33").Value = '435.77'
$ws.Range("E33").Value = '  -5.06%  '
$ws.Range("E34").Value = '  -3.09%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  -2.89%  '
$ws.Range("D37").Value = '156.85'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '19.01'
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("E40").Value = '  -5.71%  '
$ws.Range("E41").Value = '  -3.12%  '
$ws.Range("E42").Value = '  -3.44%  '
$ws.Range("E43").Value = '  -3.69%  '
$ws.Range("E44").Value = '  -4.27%  '
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("E46").Value = '  -4.99%  '
$ws.Range("D47").Value = '134.65'
$ws.Range("E47").Value = '  -4.60%  '
$ws.Range("D48").Value = '3.38'
$ws.Range("E48").Value = '  -2.83%  '
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("D50").Value = '0.487'
$ws.Range("E50").Value = '  -6.16%  '
$ws.Range("E51").Value = '  -3.12%  '
